# DFT_Energies.xlsx edit script
# Generalizes the energy-unit conversion from kcal/mol to kJ/mol, updates a
# couple of labels, tweaks a handful of data points, and moves the active
# selection — matching the target commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relabel of the excited-state manifold name ---
$ws.Range("O3").Value = '[$\bf{A}$]S$_{1-6}$'

# --- Column header: unit relabel (kcal/mol -> kJ/mol) ---
$ws.Range("G1").Value = "Relative Energy (kJ/mol)"

# --- Row 2: baseline relative energy, now in kJ/mol (x2625.5 instead of x627.5095) ---
$ws.Range("G2").Formula = "=(D2-`$D`$2)*2625.5"

# --- Row 3: literal numbers become kJ/mol-converted formulas ---
$ws.Range("G3").Formula = "=70*4.184"
$ws.Range("K3").Formula = "=7*4.184"
$ws.Range("L3").Value = 0.1

# --- Row 4 ---
$ws.Range("G4").Formula = "=(D4-`$D`$2)*2625.5"
$ws.Range("K4").Formula = "=-7*4.184"

# --- Row 5 ---
$ws.Range("G5").Formula = "=G4+(59.277*4.184)"

# --- Row 6: formula text unchanged, recalculates automatically ---

# --- Row 7 ---
$ws.Range("G7").Formula = "=(D7-`$D`$2)*2625.5"
$ws.Range("K7").Formula = "=-4*4.184"

# --- Row 8 ---
$ws.Range("G8").Formula = "=(D8-`$D`$2)*2625.5"
$ws.Range("K8").Formula = "=13*4.184"

# --- Row 9 ---
$ws.Range("G9").Formula = "=(D9-`$D`$2)*2625.5"
$ws.Range("K9").Formula = "=-16*4.184"

# --- Row 10 ---
$ws.Range("G10").Formula = "=((D10+D15)-D9)*2625.5+G9"

# --- Row 11: formula text unchanged (=G10), recalculates automatically ---

# --- Row 12 ---
$ws.Range("G12").Formula = "=(C12-C11)*2625.5+G11"
$ws.Range("K12").Formula = "=-1*4.184"

# --- Row 13 ---
$ws.Range("G13").Formula = "=((D13+D14)-D10)*2625.5+G10"
$ws.Range("K13").Formula = "=-17*4.184"

# --- Row 14: formula text unchanged (=G13), recalculates automatically ---
$ws.Range("K14").Formula = "=-1*4.184"

# --- Row 15 ---
$ws.Range("G15").Formula = "=((D15+D14)-D10)*2625.5+G10"
$ws.Range("K15").Formula = "=0*4.184"

# --- Active selection moves from L13 to K8 ---
$ws.Range("K8").Select()
